$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 21
$ws.Range("H21").Value = 2000
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 23
$ws.Range("H23").Value = 2000
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

# Row 47
$ws.Range("H47").Value = 20700
$ws.Range("J47").Value = 20700
$ws.Range("L47").Value = 20700
$ws.Range("N47").Value = -22644

# Row 64
$ws.Range("H64").Value = 2750480.8
$ws.Range("I64").Value = 5497501.5
$ws.Range("J64").Value = 3460.0715
$ws.Range("K64").Value = 5497501.5
$ws.Range("L64").Value = 3460.0715
$ws.Range("M64").Value = -5497253.5
$ws.Range("N64").Value = -3956.0715

# Row 67
$ws.Range("H67").Value = 2750480.8
$ws.Range("I67").Value = 5497501.5
$ws.Range("J67").Value = 3460.0715
$ws.Range("K67").Value = 5497501.5
$ws.Range("L67").Value = 3460.0715
$ws.Range("M67").Value = -5496643.5
$ws.Range("N67").Value = -5176.0715

# Row 93
$ws.Range("H93").Value = 36967.2
$ws.Range("J93").Value = 36967.2
$ws.Range("L93").Value = 36967.2
$ws.Range("N93").Value = -41959.2

# Row 100
$ws.Range("H100").Value = 2163.2
$ws.Range("I100").Value = 1800
$ws.Range("J100").Value = 2254
$ws.Range("K100").Value = 1800
$ws.Range("L100").Value = 2254
$ws.Range("N100").Value = -3336
$ws.Range("M100").Value = -1259

# Row 137
$ws.Range("H137").Value = 33337500
$ws.Range("I137").Value = 4999
$ws.Range("J137").Value = 40004000
$ws.Range("K137").Value = 14997
$ws.Range("L137").Value = 120012000
$ws.Range("M137").Value = -12447
$ws.Range("N137").Value = -120017100

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9418.08
$ws.Range("I32").Value = 8631.524
$ws.Range("J32").Value = 13001.277
$ws.Range("K32").Value = 8631.524
$ws.Range("L32").Value = 13001.277
$ws.Range("M32").Value = -8344.524
$ws.Range("N32").Value = -13575.277

# Row 34
$ws.Range("H34").Value = 12216.8
$ws.Range("I34").Value = 5000
$ws.Range("J34").Value = 14021
$ws.Range("K34").Value = 5000
$ws.Range("L34").Value = 14021
$ws.Range("M34").Value = -4729
$ws.Range("N34").Value = -14563

# Row 44
$ws.Range("H44").Value = 28624.125
$ws.Range("J44").Value = 28624.125
$ws.Range("L44").Value = 28624.125
$ws.Range("N44").Value = -29600.125

# Row 55
$ws.Range("H55").Value = 34111.8
$ws.Range("J55").Value = 34111.8
$ws.Range("L55").Value = 34111.8
$ws.Range("N55").Value = -34741.8

# Row 97
$ws.Range("H97").Value = 7082.625
$ws.Range("I97").Value = 8810
$ws.Range("J97").Value = 1900.5
$ws.Range("K97").Value = 8810
$ws.Range("L97").Value = 1900.5
$ws.Range("M97").Value = -8314
$ws.Range("N97").Value = -2892.5

$ws = $wb.Worksheets.Item("BSM")
# Row 36
$ws.Range("H36").Value = 1586.8334
$ws.Range("I36").Value = 1586.8334
$ws.Range("K36").Value = 1586.8334
$ws.Range("M36").Value = -1052.8334

# Row 39
$ws.Range("H39").Value = 15762
$ws.Range("J39").Value = 20000
$ws.Range("L39").Value = 20000
$ws.Range("N39").Value = -20778

# Row 56
$ws.Range("H56").Value = 44971.668
$ws.Range("J56").Value = 44971.668
$ws.Range("L56").Value = 44971.668
$ws.Range("N56").Value = -46449.668

$ws = $wb.Worksheets.Item("CRP")
# Row 8
$ws.Range("H8").Value = 500
$ws.Range("I8").Value = 0
$ws.Range("J8").Value = 500
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 500
$ws.Range("N8").Value = -780
$ws.Range("M8").ClearContents()

# Row 62
$ws.Range("H62").Value = 2485
$ws.Range("I62").Value = 2186.5625
$ws.Range("J62").Value = 3440
$ws.Range("K62").Value = 2186.5625
$ws.Range("L62").Value = 3440
$ws.Range("M62").Value = -1562.5625
$ws.Range("N62").Value = -4688

# Row 65
$ws.Range("H65").Value = 2485
$ws.Range("I65").Value = 2186.5625
$ws.Range("J65").Value = 3440
$ws.Range("K65").Value = 10932.8125
$ws.Range("L65").Value = 17200
$ws.Range("M65").Value = -7812.8125
$ws.Range("N65").Value = -23440

# Row 134
$ws.Range("H134").Value = 746287.1
$ws.Range("I134").Value = 2166.52
$ws.Range("J134").Value = 3403860.5
$ws.Range("K134").Value = 6499.559999999999
$ws.Range("L134").Value = 10211581.5
$ws.Range("M134").Value = -3964.559999999999
$ws.Range("N134").Value = -10216651.5

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 508.33334
$ws.Range("I39").Value = 100
$ws.Range("J39").Value = 545.4545
$ws.Range("K39").Value = 300
$ws.Range("L39").Value = 1636.3635
$ws.Range("N39").Value = -2224.3635
$ws.Range("M39").Value = -6

# Row 109
$ws.Range("H109").Value = 4019.2307
$ws.Range("I109").Value = 1862.5
$ws.Range("J109").Value = 4411.364
$ws.Range("K109").Value = 5587.5
$ws.Range("L109").Value = 13234.092
$ws.Range("M109").Value = -4547.5
$ws.Range("N109").Value = -15314.092

# Row 131
$ws.Range("H131").Value = 1464.4
$ws.Range("I131").Value = 2038.625
$ws.Range("J131").Value = 1194.1765
$ws.Range("K131").Value = 6115.875
$ws.Range("L131").Value = 3582.5295
$ws.Range("M131").Value = -1075.875
$ws.Range("N131").Value = -13662.5295

# Row 140
$ws.Range("H140").Value = 3529.9614
$ws.Range("I140").Value = 1431.9333
$ws.Range("J140").Value = 6390.909
$ws.Range("K140").Value = 4295.7999
$ws.Range("L140").Value = 19172.727
$ws.Range("M140").Value = 884.2001
$ws.Range("N140").Value = -29532.727

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 3993.3333
$ws.Range("J5").Value = 5000
$ws.Range("L5").Value = 5000
$ws.Range("N5").Value = -5224

# Row 63
$ws.Range("H63").Value = 39400
$ws.Range("J63").Value = 39400
$ws.Range("L63").Value = 39400
$ws.Range("N63").Value = -40772

# Row 66
$ws.Range("H66").Value = 39400
$ws.Range("J66").Value = 39400
$ws.Range("L66").Value = 118200
$ws.Range("N66").Value = -125064

$ws = $wb.Worksheets.Item("LTW")
# Row 39
$ws.Range("H39").Value = 230000
$ws.Range("I39").Value = 230000
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 230000
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("M39").Value = -229540

# Row 40
$ws.Range("H40").Value = 9445
$ws.Range("I40").Value = 50000
$ws.Range("J40").Value = 4375.625
$ws.Range("K40").Value = 50000
$ws.Range("L40").Value = 4375.625
$ws.Range("M40").Value = -49864
$ws.Range("N40").Value = -4647.625

# Row 46
$ws.Range("H46").Value = 822.48834
$ws.Range("I46").Value = 587.64703
$ws.Range("J46").Value = 976.03845
$ws.Range("K46").Value = 587.64703
$ws.Range("L46").Value = 976.03845
$ws.Range("M46").Value = -399.64703
$ws.Range("N46").Value = -1352.03845

# Row 87
$ws.Range("H87").Value = 54163
$ws.Range("J87").Value = 54163
$ws.Range("L87").Value = 54163
$ws.Range("N87").Value = -56409

# Row 90
$ws.Range("H90").Value = 54163
$ws.Range("J90").Value = 54163
$ws.Range("L90").Value = 162489
$ws.Range("N90").Value = -173721

$ws = $wb.Worksheets.Item("WVR")
# Row 23
$ws.Range("H23").Value = 1428.1111
$ws.Range("I23").Value = 786
$ws.Range("J23").Value = 2230.75
$ws.Range("K23").Value = 786
$ws.Range("L23").Value = 2230.75
$ws.Range("M23").Value = -557
$ws.Range("N23").Value = -2688.75

# Row 122
$ws.Range("H122").Value = 2700.1667
$ws.Range("I122").Value = 2621.5789
$ws.Range("J122").Value = 2998.8
$ws.Range("K122").Value = 7864.736699999999
$ws.Range("L122").Value = 8996.400000000001
$ws.Range("M122").Value = -5414.736699999999
$ws.Range("N122").Value = -13896.4
